$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the arrow glyphs from column A (rows 2-9), keeping their style.
foreach ($r in 2..9) {
    $ws.Cells.Item($r, 1).Value = ""
}

# Shrink the arrow-glyph font used by column A from 18pt to 12pt.
$ws.Range("A2:A9").Font.Size = 12

# Restore the view: selecting A2:A11 scrolls it into view (clearing the old
# topLeftCell="A7" scroll position) and updates the active cell/selection.
[void]$ws.Range("A2:A11").Select()
